$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# in_service (column E) becomes TRUE for rows 10-15
$ws.Range("E10:E15").Value = $true
